# Auto-generated script to apply cryptos list price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.167.60'
$ws.Range('E2').Value = '  -1.08%  '

$ws.Range('D3').Value = '2.642.26'
$ws.Range('E3').Value = '  -1.03%  '

$ws.Range('E4').Value = '  +0.01%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '597.21'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.80%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '156.01'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.23%  '

$ws.Range('E7').Value = '  +0.03%  '

$ws.Range('E8').Value = '  -1.13%  '

$ws.Range('E9').Value = '  +0.85%  '

$ws.Range('E10').Value = '  -1.44%  '

$ws.Range('E11').Value = '  -0.03%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.351'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.24%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '27.98'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.30%  '

$ws.Range('E14').Value = '  -0.04%  '

$ws.Range('D15').Value = '3.124.14'
$ws.Range('E15').Value = '  -0.73%  '

$ws.Range('D16').Value = '68.099.07'
$ws.Range('E16').Value = '  -1.01%  '

$ws.Range('D17').Value = '2.650.20'
$ws.Range('E17').Value = '  -0.69%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.39'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.36%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '362.99'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.26%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.35'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.43%  '

$ws.Range('E21').Value = '  +2.69%  '

$ws.Range('E23').Value = '  -3.33%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '75.12'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.40%  '

$ws.Range('E25').Value = '  -0.37%  '

$ws.Range('E26').Value = '  -4.37%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.07'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +7.11%  '

$ws.Range('E29').Value = '  -1.84%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '554.04'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -5.35%  '

$ws.Range('E31').Value = '  +0.20%  '

$ws.Range('E32').Value = '  -1.60%  '

$ws.Range('E33').Value = '  -1.20%  '

$ws.Range('E34').Value = '  +0.00%  '

$ws.Range('E35').Value = '  -2.49%  '

$ws.Range('E36').Value = '  -0.62%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '160.58'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.10%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.55'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.86%  '

$ws.Range('E39').Value = '  +0.49%  '

$ws.Range('E40').Value = '  -3.87%  '

$ws.Range('E41').Value = '  -1.79%  '

$ws.Range('E42').Value = '  +3.15%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '17.80'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.18%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.60'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.82%  '

$ws.Range('E45').Value = '  +0.03%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '158.92'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.11%  '

$ws.Range('E47').Value = '  -0.64%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '22.07'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.48%  '

$ws.Range('E49').Value = '  -2.72%  '

$ws.Range('E50').Value = '  -0.24%  '

$ws.Range('E51').Value = '  -1.05%  '
